# Update crypto price/volume figures scraped on 2023-08-12.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.432.12"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "'1.848.45"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'240.82"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.07672"
$ws.Range("E8").Value = "  +1.92%  "
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").Value = "  +1.70%  "
$ws.Range("D11").Value = "'0.07748"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "'1.846.77"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "'5.033"
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").Value = "'0.6808"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("E15").Value = "  +4.00%  "
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").Value = "'6.176"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "'29.448.57"
$ws.Range("D19").Value = "'228.26"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D22").Value = "'7.416"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'158.02"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").Value = "'0.1374"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("D26").Value = "'8.409"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("E28").Value = "  +5.75%  "
$ws.Range("D29").Value = "'1.463"
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("D30").Value = "'0.05669"
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("D31").Value = "'4.122"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").Value = "'4.028"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").Value = "'1.844"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").Value = "'1.162"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("D35").Value = "'0.7025"
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("D36").Value = "'2.590"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").Value = "'1.226.48"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("D40").Value = "'6.543"
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("D41").Value = "'0.9048"
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("D43").Value = "'2.001.86"
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("D44").Value = "'101.80"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "'66.04"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").Value = "'0.00000000121"
$ws.Range("E46").Value = "  +2.35%  "
$ws.Range("D47").Value = "'7.163"
$ws.Range("E47").Value = "  +1.49%  "
$ws.Range("D48").Value = "'0.4018"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("D49").Value = "'0.1156"
$ws.Range("E49").Value = "  +3.38%  "
$ws.Range("D50").Value = "'9.010"
$ws.Range("E50").Value = "  +1.59%  "
$ws.Range("E51").Value = "  +0.52%  "
